$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear avatar (column F) data cells that referenced the removed shared strings
$ws.Range("F5").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("F10").ClearContents()

# Set column widths for E (email) and F (avatar) to match bestFit sizing
$ws.Columns.Item(5).ColumnWidth = 22.5
$ws.Columns.Item(6).ColumnWidth = 50.39

# Update the active selection on the sheet view
$ws.Range("F12").Select()
